$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.674.63"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.886.69"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -1.19%  "
$ws.Range("D5").Value = "'313.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "'0.4859"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("D8").Value = "'0.3784"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'0.9179"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "'20.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07675"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.918.40"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "'5.460"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "'6.579"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'90.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'0.000008792"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "27.724.38"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'14.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "'5.111"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "2.164.94"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").Value = "'1.901"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "'153.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "'18.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "'2.107"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.89%  "
$ws.Range("D29").Value = "'115.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'4.892"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "'0.08935"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'3.146"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.59%  "
$ws.Range("D33").Value = "'1.218"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'0.7564"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "'0.02037"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'2.528"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("D38").Value = "'1.089"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("D39").Value = "'0.05241"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").Value = "'0.5437"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").Value = "'2.969"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'6.936"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "'8.308"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "'109.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("D46").Value = "'10.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'0.4770"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'1.625"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "'0.06053"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
